$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Health pickup sound implemented: asset finalized (dropped the "(x2)"
# duplicate-variant note) and the row marked Completed.
$ws.Range("D22").Value = "Health increase SFX"
$ws.Range("E22").Value = "Completed"

# Match the author's final cursor/selection position.
$ws.Range("E21").Select()
